$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.488.27'
$ws.Range('D3').Value = '2.128.77'
$ws.Range('E3').Value = '  +1.53%  '
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').Value = '348.25'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').Value = '1.004'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').Value = '0.5231'
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('D8').Value = '0.4495'
$ws.Range('E8').Value = '  +0.74%  '
$ws.Range('D9').Value = '54.18'
$ws.Range('E9').Value = '  +4.12%  '
$ws.Range('D10').Value = '0.09407'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').Value = '1.187'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').Value = '25.51'
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').Value = '8.749'
$ws.Range('E13').Value = '  +8.57%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '6.981'
$ws.Range('E14').Value = '  +3.36%  '
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '2.083.24'
$ws.Range('E15').Value = '  -0.75%  '
$ws.Range('D16').Value = '103.36'
$ws.Range('E16').Value = '  +4.09%  '
$ws.Range('D17').Value = '0.00001173'
$ws.Range('E17').Value = '  +0.88%  '
$ws.Range('E18').Value = '  -0.40%  '
$ws.Range('D19').Value = '21.61'
$ws.Range('E19').Value = '  +4.89%  '
$ws.Range('D20').Value = '0.06707'
$ws.Range('E20').Value = '  -0.03%  '
$ws.Range('D21').Value = '6.343'
$ws.Range('E21').Value = '  +2.58%  '
$ws.Range('D22').Value = '1.004'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').Value = '30.483.92'
$ws.Range('E23').Value = '  +1.96%  '
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').Value = '2.338'
$ws.Range('E25').Value = '  +0.89%  '
$ws.Range('D26').Value = '2.374.85'
$ws.Range('E26').Value = '  +1.30%  '
$ws.Range('E27').Value = '  +1.30%  '
$ws.Range('D28').Value = '2.561'
$ws.Range('E28').Value = '  +0.75%  '
$ws.Range('D29').Value = '163.86'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('D30').Value = '134.84'
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('D31').Value = '1.173'
$ws.Range('E31').Value = '  +0.84%  '
$ws.Range('D32').Value = '1.809'
$ws.Range('E32').Value = '  +11.58%  '
$ws.Range('E33').Value = '  +0.29%  '
$ws.Range('D34').Value = '6.965'
$ws.Range('E34').Value = '  +13.26%  '
$ws.Range('E35').Value = '  +1.57%  '
$ws.Range('D36').Value = '3.957'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = '10.74'
$ws.Range('E37').Value = '  +5.76%  '
$ws.Range('D38').Value = '0.02654'
$ws.Range('E38').Value = '  +3.33%  '
$ws.Range('D39').Value = '0.06887'
$ws.Range('E39').Value = '  +1.97%  '
$ws.Range('D40').Value = '0.7179'
$ws.Range('E40').Value = '  +4.02%  '
$ws.Range('E41').Value = '  +2.65%  '
$ws.Range('E42').Value = '  -0.78%  '
$ws.Range('D43').Value = '1.338'
$ws.Range('E43').Value = '  +2.98%  '
$ws.Range('D44').Value = '0.6972'
$ws.Range('E44').Value = '  +4.97%  '
$ws.Range('E45').Value = '  +4.85%  '
$ws.Range('D46').Value = '2.404'
$ws.Range('E46').Value = '  +5.03%  '
$ws.Range('D47').Value = '1.004'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('D48').Value = '1.279'
$ws.Range('E48').Value = '  +9.48%  '
$ws.Range('D49').Value = '3.635'
$ws.Range('E49').Value = '  -0.21%  '
$ws.Range('D50').Value = '0.00000000347'
$ws.Range('E50').Value = '  +1.64%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '1.232'
$ws.Range('E51').Value = '  +1.09%  '
